$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.831.72'
$ws.Range('E2').Value = '  -0.74%  '
$ws.Range('D3').Value = '2.351.56'
$ws.Range('E3').Value = '  -0.39%  '
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('E5').Value = '  -1.46%  '
$ws.Range('D6').Value = "'240.38"
$ws.Range('E6').Value = '  -1.47%  '
$ws.Range('D7').Value = "'73.17"
$ws.Range('E7').Value = '  -1.69%  '
$ws.Range('E8').Value = '  -0.03%  '
$ws.Range('E9').Value = '  +2.62%  '
$ws.Range('D10').Value = "'0.100"
$ws.Range('E10').Value = '  -2.61%  '
$ws.Range('D11').Value = "'60.25"
$ws.Range('E11').Value = '  +4.44%  '
$ws.Range('D12').Value = "'32.96"
$ws.Range('E12').Value = '  +3.27%  '
$ws.Range('D13').Value = "'7.32"
$ws.Range('E13').Value = '  -2.17%  '
$ws.Range('E14').Value = '  -0.38%  '
$ws.Range('D15').Value = '2.700.25'
$ws.Range('E15').Value = '  -0.38%  '
$ws.Range('D16').Value = "'16.37"
$ws.Range('E16').Value = '  -3.20%  '
$ws.Range('D17').Value = "'0.904"
$ws.Range('E17').Value = '  -1.24%  '
$ws.Range('D18').Value = '2.342.68'
$ws.Range('E18').Value = '  -0.89%  '
$ws.Range('D19').Value = '43.755.83'
$ws.Range('E19').Value = '  -1.49%  '
$ws.Range('E20').Value = '  -2.13%  '
$ws.Range('D21').Value = "'6.69"
$ws.Range('E21').Value = '  -1.07%  '
$ws.Range('D22').Value = "'77.25"
$ws.Range('E22').Value = '  -1.57%  '
$ws.Range('D23').Value = "'254.03"
$ws.Range('E23').Value = '  -0.89%  '
$ws.Range('E24').Value = '  +17.90%  '
$ws.Range('E25').Value = '  -0.03%  '
$ws.Range('D26').Value = "'3.73"
$ws.Range('E26').Value = '  -0.72%  '
$ws.Range('E27').Value = '  -3.50%  '
$ws.Range('B28').Value = 'Cosmos'
$ws.Range('C28').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D28').Value = "'10.55"
$ws.Range('E28').Value = '  -1.82%  '
$ws.Range('B29').Value = 'Toncoin'
$ws.Range('C29').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D29').Value = "'2.27"
$ws.Range('E29').Value = '  -1.56%  '
$ws.Range('B30').Value = 'EthereumClassic'
$ws.Range('C30').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D30').Value = "'22.64"
$ws.Range('E30').Value = '  +0.17%  '
$ws.Range('B31').Value = 'Monero'
$ws.Range('C31').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D31').Value = "'177.43"
$ws.Range('E31').Value = '  +1.48%  '
$ws.Range('E32').Value = '  -1.36%  '
$ws.Range('E33').Value = '  +1.92%  '
$ws.Range('D34').Value = "'0.0757"
$ws.Range('E34').Value = '  -0.13%  '
$ws.Range('D36').Value = "'5.43"
$ws.Range('E36').Value = '  +1.25%  '
$ws.Range('E37').Value = '  -2.13%  '
$ws.Range('E38').Value = '  -3.86%  '
$ws.Range('E39').Value = '  -3.63%  '
$ws.Range('D40').Value = "'0.0278"
$ws.Range('E40').Value = '  +1.11%  '
$ws.Range('D41').Value = "'68.55"
$ws.Range('E41').Value = '  +29.24%  '
$ws.Range('E42').Value = '  +11.53%  '
$ws.Range('B43').Value = 'FTXToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D43').Value = "'5.01"
$ws.Range('E43').Value = '  +11.54%  '
$ws.Range('B44').Value = 'FraxShare'
$ws.Range('C44').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D44').Value = "'9.14"
$ws.Range('E44').Value = '  +1.36%  '
$ws.Range('B45').Value = 'InjectiveProtocol'
$ws.Range('C45').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D45').Value = "'19.10"
$ws.Range('E45').Value = '  -1.77%  '
$ws.Range('E46').Value = '  +3.10%  '
$ws.Range('E47').Value = '  -0.74%  '
$ws.Range('E48').Value = '  -0.92%  '
$ws.Range('E49').Value = '  +0.14%  '
$ws.Range('B50').Value = 'Aave'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D50').Value = "'99.05"
$ws.Range('E50').Value = '  -2.61%  '
$ws.Range('B51').Value = 'ARBITRUM'
$ws.Range('C51').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D51').Value = "'1.16"
$ws.Range('E51').Value = '  -1.69%  '
